{"js": "// Insert four new bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n// section, right after the \"Platform impact\" bullet and before the\n// \"TECHNICAL SKILLS\" heading.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph: the \"Platform impact\" bullet.\nconst anchorText = \"Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for 'Platform impact' bullet.\");\n}\n\n// New bullet lines to insert, in order, after the anchor paragraph.\nconst newLines = [\n  \"\u2022 Real-time collaboration at national scale\",\n  \"\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n  \"\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n  \"\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"\n];\n\n// Insert each line right after the anchor, one after another, so the\n// final order matches newLines (insert \"After\" anchor repeatedly while\n// always using the anchor as the reference keeps pushing each new\n// paragraph directly below the anchor; instead, chain from the\n// previously inserted paragraph to preserve ordering).\nlet previous = anchor;\nfor (let i = 0; i < newLines.length; i++) {\n  previous = previous.insertParagraph(newLines[i], \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert four new bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n# section, right after the \"Platform impact\" bullet and before the\n# \"TECHNICAL SKILLS\" heading.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"Platform impact: ...\") using Find so we\n# don't depend on a hard-coded paragraph index.\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find anchor paragraph for 'Platform impact' bullet.\"\n}\n\n# End of the matched text == end of the anchor paragraph's content\n# (just before its paragraph mark).\n$pos = $searchRange.End\n\n$newLines = @(\n    \"\u2022 Real-time collaboration at national scale\",\n    \"\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%\",\n    \"\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    \"\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"\n)\n\nforeach ($line in $newLines) {\n    # Insert a new paragraph break right after the current position.\n    $breakRange = $d.Range($pos, $pos)\n    $breakRange.InsertParagraphAfter()\n\n    # The new (empty) paragraph now starts one character past $pos (the\n    # paragraph mark that was just inserted). Fill it with the bullet text.\n    $pos = $pos + 1\n    $textRange = $d.Range($pos, $pos)\n    $textRange.InsertAfter($line)\n\n    # Advance past the text we just inserted so the next paragraph break\n    # is added after this new bullet.\n    $pos = $pos + $line.Length\n}\n"}
